$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $origStyle = $rng.Style
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = $origStyle
}

Set-TextValue "D2" '61.073.93'
Set-TextValue "E2" '  -0.06%  '

Set-TextValue "D3" '3.366.79'
Set-TextValue "E3" '  +1.64%  '

Set-TextValue "E4" '  -0.07%  '

Set-TextValue "D5" '571.52'
Set-TextValue "E5" '  +1.48%  '

Set-TextValue "D6" '136.79'
Set-TextValue "E6" '  +7.04%  '

Set-TextValue "E7" '  -0.03%  '

Set-TextValue "D8" '3.365.89'
Set-TextValue "E8" '  +1.65%  '

Set-TextValue "D9" '0.476'
Set-TextValue "E9" '  -0.29%  '

Set-TextValue "E10" '  +3.85%  '

Set-TextValue "D11" '0.123'
Set-TextValue "E11" '  +4.25%  '

Set-TextValue "D12" '0.392'
Set-TextValue "E12" '  +5.02%  '

Set-TextValue "D13" '3.940.27'
Set-TextValue "E13" '  +1.52%  '

Set-TextValue "E14" '  +2.03%  '

Set-TextValue "D15" '0.0000173'
Set-TextValue "E15" '  +3.03%  '

Set-TextValue "D16" '3.364.81'
Set-TextValue "E16" '  +1.45%  '

Set-TextValue "D17" '25.22'
Set-TextValue "E17" '  +1.77%  '

Set-TextValue "D18" '61.171.51'
Set-TextValue "E18" '  -0.07%  '

$ws.Range("B19").Value = 'Polkadot'
$ws.Range("C19").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
Set-TextValue "D19" '5.88'
Set-TextValue "E19" '  +4.73%  '

$ws.Range("B20").Value = 'Chainlink'
$ws.Range("C20").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
Set-TextValue "D20" '13.89'
Set-TextValue "E20" '  +3.91%  '

Set-TextValue "D21" '9.35'
Set-TextValue "E21" '  +4.19%  '

Set-TextValue "D22" '378.33'
Set-TextValue "E22" '  +7.48%  '

Set-TextValue "D23" '0.567'
Set-TextValue "E23" '  +2.40%  '

Set-TextValue "D24" '3.501.06'
Set-TextValue "E24" '  +1.57%  '

Set-TextValue "E25" '  -0.02%  '

Set-TextValue "D26" '70.59'
Set-TextValue "E26" '  +2.07%  '

Set-TextValue "E27" '  +10.92%  '

Set-TextValue "D28" '1.65'
Set-TextValue "E28" '  +17.94%  '

Set-TextValue "D29" '7.75'
Set-TextValue "E29" '  +8.24%  '

Set-TextValue "E30" '  +0.19%  '

Set-TextValue "D31" '8.18'
Set-TextValue "E31" '  +4.14%  '

Set-TextValue "E32" '  +4.85%  '

Set-TextValue "E33" '  +1.19%  '

Set-TextValue "E34" '  -0.08%  '

Set-TextValue "D35" '3.397.80'
Set-TextValue "E35" '  +1.68%  '

Set-TextValue "D36" '23.36'
Set-TextValue "E36" '  +3.74%  '

Set-TextValue "D37" '5.59'
Set-TextValue "E37" '  +7.21%  '

Set-TextValue "D38" '7.02'
Set-TextValue "E38" '  +3.73%  '

Set-TextValue "E39" '  +5.10%  '

Set-TextValue "D40" '161.20'
Set-TextValue "E40" '  +0.27%  '

Set-TextValue "D41" '0.0790'
Set-TextValue "E41" '  +4.33%  '

Set-TextValue "E42" '  -0.08%  '

$ws.Range("B43").Value = 'Filecoin'
$ws.Range("C43").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
Set-TextValue "D43" '4.42'
Set-TextValue "E43" '  +2.19%  '

$ws.Range("B44").Value = 'Stacks'
$ws.Range("C44").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
Set-TextValue "D44" '1.70'
Set-TextValue "E44" '  +8.61%  '

Set-TextValue "D45" '41.43'
Set-TextValue "E45" '  +1.07%  '

$ws.Range("B46").Value = 'Mantle'
$ws.Range("C46").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
Set-TextValue "D46" '0.758'
Set-TextValue "E46" '  +2.08%  '

$ws.Range("B47").Value = 'ONDO'
$ws.Range("C47").Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
Set-TextValue "D47" '1.19'
Set-TextValue "E47" '  +7.46%  '

Set-TextValue "D48" '23.09'
Set-TextValue "E48" '  +3.77%  '

Set-TextValue "E49" '  +4.18%  '

Set-TextValue "D50" '22.85'
Set-TextValue "E50" '  +9.49%  '

Set-TextValue "D51" '2.317.90'
Set-TextValue "E51" '  +6.74%  '
